$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ------------------------------------------------------------------
# Reset the working area: delete the rows outright so that stale
# content, styles and custom row-heights are fully cleared (a plain
# ClearContents/Clear leaves old row-height metadata behind).
# ------------------------------------------------------------------
$ws.Range("A1:K20").EntireRow.Delete()

# ------------------------------------------------------------------
# Re-populate the table with the new layout
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Table"
$ws.Range("B1").Value = "Traces"
$ws.Range("C1").Value = "Traces"
$ws.Range("D1").Value = "Dependencies HTTP"
$ws.Range("E1").Value = "Request(Http)"
$ws.Range("F1").Value = "Exceptions"
$ws.Range("G1").Value = "Exceptions"
$ws.Range("H1").Value = "Metric"
$ws.Range("I1").Value = "Event"
$ws.Range("J1").Value = "Dependencies HTTP"
$ws.Range("K1").Value = "Request(Http)"

$ws.Range("A2").Value = "Integration type"
$ws.Range("B2").Value = "Serilog-TelemetryClient"
$ws.Range("C2").Value = "AppService"
$ws.Range("D2").Value = "Serilog-TelemetryClient"
$ws.Range("E2").Value = "AppService"
$ws.Range("F2").Value = "Serilog-TelemetryClient"
$ws.Range("G2").Value = "AppService"
$ws.Range("H2").Value = "Serilog-TelemetryClient"
$ws.Range("I2").Value = "Serilog-TelemetryClient"
$ws.Range("J2").Value = "AppService"
$ws.Range("K2").Value = "Serilog-TelemetryClient"

$ws.Range("A3").Value = "Present"
$ws.Range("B3").Value = "Yes"
$ws.Range("C3").Value = "Yes"
$ws.Range("D3").Value = "Yes"
$ws.Range("E3").Value = "Yes"
$ws.Range("F3").Value = "TBD"
$ws.Range("H3").Value = "TBD"
$ws.Range("I3").Value = "TBD"
$ws.Range("J3").Value = "No"
$ws.Range("K3").Value = "No"

$ws.Range("A4").Value = "Custom Properties"
$ws.Range("B4").Value = "Yes"
$ws.Range("C4").Value = "No"
$ws.Range("D4").Value = "No"
$ws.Range("J4").Value = "NA"
$ws.Range("K4").Value = "NA"

$ws.Range("A5").Value = "Cloud application identifier"
$ws.Range("B5").Value = "cloud_RoleName"
$ws.Range("C5").Value = "cloud_RoleName"
$ws.Range("D5").Value = "cloud_RoleName"
$ws.Range("E5").Value = "cloud_RoleName"
$ws.Range("J5").Value = "NA"
$ws.Range("K5").Value = "NA"

$ws.Range("A6").Value = "Cloud Instance Identifier"
$ws.Range("B6").Value = "customDimensions.HostInstanceId"
$ws.Range("C6").Value = "customDimensions.HostInstanceId"
$ws.Range("E6").Value = "customDimensions.HostInstanceId"
$ws.Range("J6").Value = "NA"
$ws.Range("K6").Value = "NA"

$ws.Range("A7").Value = "Executing identifier"
$ws.Range("B7").Value = "customDimensions.InvocationId"
$ws.Range("C7").Value = "No (As all the logs have come to the console)`nNeed to check about "
$ws.Range("E7").Value = "customDimensions.InvocationId"
$ws.Range("J7").Value = "NA"
$ws.Range("K7").Value = "NA"

$ws.Range("A8").Value = "Distributed transaction identifier"
$ws.Range("B8").Value = "TBD"
$ws.Range("J8").Value = "NA"
$ws.Range("K8").Value = "NA"

$ws.Range("A9").Value = "Correlatable Fields"
$ws.Range("J9").Value = "NA"
$ws.Range("K9").Value = "NA"

$ws.Range("A10").Value = "Sampling"
$ws.Range("D10").Value = "Dependencies are getting sampled?"
$ws.Range("J10").Value = "NA"
$ws.Range("K10").Value = "NA"

$ws.Range("A11").Value = "Remarks"
$ws.Range("D11").Value = "Manually tracked for HTTP`nThe cloud role name has to be set"
$ws.Range("J11").Value = "NA"
$ws.Range("K11").Value = "NA"

$ws.Range("A12").Value = "UserCorrelation"
$ws.Range("B12").Value = "TBD"

$ws.Range("A13").Value = "FileCorrelation/Module Correlation"
$ws.Range("B13").Value = "customDimensions.Category"
$ws.Range("C13").Value = "NA"
$ws.Range("D13").Value = "Relook needed (Using deprecated API)"

$ws.Range("A14").Value = "Per log property for traces"

$ws.Range("A15").Value = "appID?"

# ------------------------------------------------------------------
# Wrap-text formatting on the same cells that had it before
# ------------------------------------------------------------------
$ws.Range("C7").WrapText = $true
$ws.Range("G8").WrapText = $true
$ws.Range("C9").WrapText = $true
$ws.Range("D9").WrapText = $true
$ws.Range("C11").WrapText = $true
$ws.Range("D11").WrapText = $true
$ws.Range("D13").WrapText = $true

# ------------------------------------------------------------------
# Row heights for the two wrapped, two-line remark rows
# ------------------------------------------------------------------
$ws.Rows.Item(7).RowHeight = 28.5
$ws.Rows.Item(11).RowHeight = 28.5

# ------------------------------------------------------------------
# Column widths
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 43.5
$ws.Columns.Item(3).ColumnWidth = 35.67
$ws.Columns.Item(6).ColumnWidth = 18.5
$ws.Columns.Item(7).ColumnWidth = 13.17
$ws.Columns.Item(8).ColumnWidth = 21.5
$ws.Columns.Item(9).ColumnWidth = 19.17

# ------------------------------------------------------------------
# Freeze panes: split after column A, keep row 1 unsplit, select D25
# ------------------------------------------------------------------
$ws.Activate()
$ws.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D25").Select()

# ------------------------------------------------------------------
# Page setup - portrait orientation
# ------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
